# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)
#
# Updates the "Enterprises density (per 1000 people)" statistics on the
# Azerbaijan MSME summary sheet with more precise figures, for both the
# "Statistical Institution" table (row 11) and the "SME Associations"
# table (row 28). The values are kept as text (matching how the rest of
# the sheet stores its numeric-looking figures as shared strings) by
# prefixing the literal with an apostrophe; ClearFormats() afterwards
# drops the resulting quote-prefix / text-format styling so the cell
# keeps its original (default) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Source Type: Statistical Institution (Most Widely Used)
# Enterprises density (per 1000 people): Micro / SMEs / MSMEs
$ws.Range("B11").Value = "'26.71"
$ws.Range("B11").ClearFormats()
$ws.Range("C11").Value = "'1.47"
$ws.Range("C11").ClearFormats()
$ws.Range("D11").Value = "'28.17"
$ws.Range("D11").ClearFormats()

# Source Type: SME Associations
# Enterprises density (per 1000 people): Micro / SMEs / MSMEs
$ws.Range("B28").Value = "'21.27"
$ws.Range("B28").ClearFormats()
$ws.Range("C28").Value = "'8.15"
$ws.Range("C28").ClearFormats()
$ws.Range("D28").Value = "'29.42"
$ws.Range("D28").ClearFormats()
